# Applies the cryptos price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.973.83"
$ws.Cells.Item(2, 5).Value = "  +2.08%  "
$ws.Cells.Item(3, 4).Value = "1.903.83"
$ws.Cells.Item(3, 5).Value = "  +1.85%  "
$ws.Cells.Item(4, 4).Value = "'1.001"
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).Value = "'333.02"
$ws.Cells.Item(5, 5).Value = "  -1.78%  "
$ws.Cells.Item(6, 5).Value = "  +0.04%  "
$ws.Cells.Item(7, 4).Value = "'0.4627"
$ws.Cells.Item(7, 5).Value = "  -1.45%  "
$ws.Cells.Item(8, 4).Value = "'0.4083"
$ws.Cells.Item(8, 5).Value = "  +2.75%  "
$ws.Cells.Item(9, 4).Value = "'47.85"
$ws.Cells.Item(9, 5).Value = "  +0.53%  "
$ws.Cells.Item(10, 4).Value = "'0.08028"
$ws.Cells.Item(10, 5).Value = "  +0.10%  "
$ws.Cells.Item(11, 4).Value = "'1.004"
$ws.Cells.Item(11, 5).Value = "  +0.33%  "
$ws.Cells.Item(12, 4).Value = "'21.73"
$ws.Cells.Item(12, 5).Value = "  -1.02%  "
$ws.Cells.Item(13, 4).Value = "1.913.65"
$ws.Cells.Item(13, 5).Value = "  +1.63%  "
$ws.Cells.Item(14, 4).Value = "'5.939"
$ws.Cells.Item(15, 4).Value = "'7.077"
$ws.Cells.Item(16, 4).Value = "'1.001"
$ws.Cells.Item(16, 5).Value = "  -0.14%  "
$ws.Cells.Item(17, 4).Value = "'88.87"
$ws.Cells.Item(17, 5).Value = "  -2.63%  "
$ws.Cells.Item(18, 4).Value = "'0.00001033"
$ws.Cells.Item(18, 5).Value = "  -0.88%  "
$ws.Cells.Item(19, 4).Value = "'0.06573"
$ws.Cells.Item(19, 5).Value = "  -0.63%  "
$ws.Cells.Item(20, 4).Value = "'17.52"
$ws.Cells.Item(20, 5).Value = "  -0.18%  "
$ws.Cells.Item(21, 4).Value = "'1.001"
$ws.Cells.Item(21, 5).Value = "  +0.00%  "
$ws.Cells.Item(22, 4).Value = "28.992.01"
$ws.Cells.Item(22, 5).Value = "  +2.11%  "
$ws.Cells.Item(23, 4).Value = "'5.455"
$ws.Cells.Item(23, 5).Value = "  -0.36%  "
$ws.Cells.Item(24, 4).Value = "'11.29"
$ws.Cells.Item(24, 5).Value = "  +2.04%  "
$ws.Cells.Item(25, 4).Value = "'2.232"
$ws.Cells.Item(25, 5).Value = "  -1.13%  "
$ws.Cells.Item(26, 4).Value = "2.136.17"
$ws.Cells.Item(26, 5).Value = "  +1.51%  "
$ws.Cells.Item(27, 4).Value = "'157.78"
$ws.Cells.Item(27, 5).Value = "  -2.12%  "
$ws.Cells.Item(28, 4).Value = "'19.73"
$ws.Cells.Item(28, 5).Value = "  -0.18%  "
$ws.Cells.Item(29, 4).Value = "'2.098"
$ws.Cells.Item(29, 5).Value = "  -1.20%  "
$ws.Cells.Item(30, 5).Value = "  -1.93%  "
$ws.Cells.Item(31, 4).Value = "'118.81"
$ws.Cells.Item(31, 5).Value = "  -1.28%  "
$ws.Cells.Item(32, 4).Value = "'0.9787"
$ws.Cells.Item(32, 5).Value = "  +0.95%  "
$ws.Cells.Item(33, 4).Value = "'0.09406"
$ws.Cells.Item(33, 5).Value = "  -1.07%  "
$ws.Cells.Item(34, 4).Value = "'1.419"
$ws.Cells.Item(34, 5).Value = "  +3.25%  "
$ws.Cells.Item(35, 4).Value = "'3.586"
$ws.Cells.Item(35, 5).Value = "  -0.29%  "
$ws.Cells.Item(36, 4).Value = "'5.311"
$ws.Cells.Item(36, 5).Value = "  -0.73%  "
$ws.Cells.Item(37, 4).Value = "'0.06087"
$ws.Cells.Item(38, 4).Value = "'0.02240"
$ws.Cells.Item(38, 5).Value = "  -0.58%  "
$ws.Cells.Item(39, 4).Value = "'8.388"
$ws.Cells.Item(39, 5).Value = "  +0.18%  "
$ws.Cells.Item(40, 4).Value = "'1.170"
$ws.Cells.Item(40, 5).Value = "  -1.30%  "
$ws.Cells.Item(41, 4).Value = "'0.5812"
$ws.Cells.Item(41, 5).Value = "  -2.32%  "
$ws.Cells.Item(42, 4).Value = "'1.000"
$ws.Cells.Item(42, 5).Value = "  -0.02%  "
$ws.Cells.Item(43, 4).Value = "'10.17"
$ws.Cells.Item(43, 5).Value = "  -1.72%  "
$ws.Cells.Item(44, 4).Value = "'0.1821"
$ws.Cells.Item(44, 5).Value = "  -2.80%  "
$ws.Cells.Item(45, 4).Value = "'1.248"
$ws.Cells.Item(45, 5).Value = "  -3.12%  "
$ws.Cells.Item(46, 4).Value = "'2.302"
$ws.Cells.Item(46, 5).Value = "  +11.08%  "
$ws.Cells.Item(47, 4).Value = "'12.15"
$ws.Cells.Item(47, 5).Value = "  -0.15%  "
$ws.Cells.Item(48, 4).Value = "'0.5503"
$ws.Cells.Item(49, 4).Value = "'1.908"
$ws.Cells.Item(49, 5).Value = "  -2.70%  "
$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(50, 4).Value = "'48.49"
$ws.Cells.Item(50, 5).Value = "  +25.11%  "
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).Value = "'0.07027"
$ws.Cells.Item(51, 5).Value = "  +2.31%  "
